$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header row, shifting the existing
# lookup rows down by one.
$ws.Rows(2).Insert()

# Populate the new "Not applicable" lookup entry.
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# The named range covering the table grew by one row (it was $A$1:$D$62,
# now spans the newly-inserted row too).
$wb.Names("dbo_bodypar").RefersTo = "=dbo_bodypar!`$A`$1:`$D`$63"

# Match the saved selection/active cell state.
$null = $ws.Range("B3").Select()
